$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 200
$ws.Range("I6").Value = 200
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 600
$ws.Range("L6").Value = 600
$ws.Range("M6").Value = -488
$ws.Range("N6").Value = -824
$ws.Range("H8").Value = 40.8
$ws.Range("I8").Value = 20.75
$ws.Range("J8").Value = 121
$ws.Range("K8").Value = 62.25
$ws.Range("L8").Value = 363
$ws.Range("M8").Value = 76.75
$ws.Range("N8").Value = -641
$ws.Range("H31").Value = 1075.25
$ws.Range("I31").Value = 17
$ws.Range("K31").Value = 51
$ws.Range("M31").Value = 179
$ws.Range("H86").Value = 9625.817999999999
$ws.Range("I86").Value = 10720.444
$ws.Range("J86").Value = 4700
$ws.Range("K86").Value = 10720.444
$ws.Range("L86").Value = 4700
$ws.Range("M86").Value = -9597.444
$ws.Range("N86").Value = -6946
$ws.Range("H89").Value = 9625.817999999999
$ws.Range("I89").Value = 10720.444
$ws.Range("J89").Value = 4700
$ws.Range("K89").Value = 53602.22
$ws.Range("L89").Value = 23500
$ws.Range("M89").Value = -47986.22
$ws.Range("N89").Value = -34732
$ws.Range("H98").Value = 2754
$ws.Range("I98").Value = 2468.9333
$ws.Range("J98").Value = 3823
$ws.Range("K98").Value = 2468.9333
$ws.Range("L98").Value = 3823
$ws.Range("M98").Value = -970.9333000000001
$ws.Range("N98").Value = -6819
$ws.Range("H112").Value = 2572.5217
$ws.Range("H122").Value = 2754
$ws.Range("I122").Value = 2468.9333
$ws.Range("J122").Value = 3823
$ws.Range("K122").Value = 7406.7999
$ws.Range("L122").Value = 11469
$ws.Range("M122").Value = -4956.7999
$ws.Range("N122").Value = -16369
$ws.Range("H135").Value = 1135.7333
$ws.Range("I135").Value = 1028
$ws.Range("K135").Value = 9252
$ws.Range("M135").Value = -6717

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 2123.8
$ws.Range("I5").Value = 109.5
$ws.Range("J5").Value = 3466.6667
$ws.Range("K5").Value = 109.5
$ws.Range("L5").Value = 3466.6667
$ws.Range("M5").Value = 3.5
$ws.Range("N5").Value = -3692.6667
$ws.Range("H86").Value = 3202.8572
$ws.Range("I86").Value = 1934.7
$ws.Range("K86").Value = 1934.7
$ws.Range("M86").Value = -811.7
$ws.Range("H89").Value = 3202.8572
$ws.Range("I89").Value = 1934.7
$ws.Range("K89").Value = 9673.5
$ws.Range("M89").Value = -4057.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 20333
$ws.Range("I23").Value = 20499.5
$ws.Range("J23").Value = 20000
$ws.Range("K23").Value = 20499.5
$ws.Range("L23").Value = 20000
$ws.Range("M23").Value = -20259.5
$ws.Range("N23").Value = -20480
$ws.Range("H27").Value = 20333
$ws.Range("I27").Value = 20499.5
$ws.Range("J27").Value = 20000
$ws.Range("K27").Value = 20499.5
$ws.Range("L27").Value = 20000
$ws.Range("M27").Value = -20307.5
$ws.Range("N27").Value = -20384
$ws.Range("H31").Value = 3224.6
$ws.Range("J31").Value = 3338
$ws.Range("L31").Value = 3338
$ws.Range("N31").Value = -3928
$ws.Range("H34").Value = 3224.6
$ws.Range("J34").Value = 3338
$ws.Range("L34").Value = 3338
$ws.Range("N34").Value = -3742
$ws.Range("H86").Value = 8166.6665
$ws.Range("J86").Value = 6000
$ws.Range("L86").Value = 6000
$ws.Range("N86").Value = -8246
$ws.Range("H89").Value = 8166.6665
$ws.Range("J89").Value = 6000
$ws.Range("L89").Value = 30000
$ws.Range("N89").Value = -41232
$ws.Range("H134").Value = 1860.4642
$ws.Range("I134").Value = 1088.5454
$ws.Range("K134").Value = 3265.6362
$ws.Range("M134").Value = -730.6361999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 506.2
$ws.Range("I2").Value = 52.6
$ws.Range("K2").Value = 315.6
$ws.Range("M2").Value = -202.6
$ws.Range("H11").Value = 3997.5
$ws.Range("J11").Value = 3997.5
$ws.Range("L11").Value = 11992.5
$ws.Range("N11").Value = -12272.5
$ws.Range("H26").Value = 262675
$ws.Range("J26").Value = 25300
$ws.Range("L26").Value = 75900
$ws.Range("N26").Value = -76476
$ws.Range("H33").Value = 25
$ws.Range("J33").Value = 25
$ws.Range("L33").Value = 150
$ws.Range("N33").Value = -716
$ws.Range("H107").Value = 798
$ws.Range("J107").Value = 815.6667
$ws.Range("L107").Value = 2447.0001
$ws.Range("N107").Value = -6287.0001
$ws.Range("H131").Value = 2433.3333
$ws.Range("J131").Value = 2433.3333
$ws.Range("L131").Value = 7299.999899999999
$ws.Range("N131").Value = -17379.9999
$ws.Range("H138").Value = 7500.1665
$ws.Range("J138").Value = 8750.25
$ws.Range("L138").Value = 26250.75
$ws.Range("N138").Value = -36530.75
$ws.Range("H141").Value = 1996.5
$ws.Range("I141").Value = 1996.5
$ws.Range("K141").Value = 5989.5
$ws.Range("M141").Value = -809.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 12000
$ws.Range("J39").Value = 12000
$ws.Range("L39").Value = 12000
$ws.Range("N39").Value = -13064
$ws.Range("H80").Value = 4219.4443
$ws.Range("J80").Value = 4219.4443
$ws.Range("L80").Value = 4219.4443
$ws.Range("N80").Value = -6215.4443
$ws.Range("H83").Value = 4219.4443
$ws.Range("J83").Value = 4219.4443
$ws.Range("L83").Value = 21097.2215
$ws.Range("N83").Value = -31081.2215
$ws.Range("H122").Value = 11368772
$ws.Range("I122").Value = 13894269
$ws.Range("J122").Value = 4034
$ws.Range("K122").Value = 41682807
$ws.Range("L122").Value = 12102
$ws.Range("M122").Value = -41680357
$ws.Range("N122").Value = -17002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3127.3572
$ws.Range("I122").Value = 3106.3845
$ws.Range("K122").Value = 9319.1535
$ws.Range("M122").Value = -6869.1535

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 3000.5
$ws.Range("I13").Value = 5001
$ws.Range("J13").Value = 1000
$ws.Range("K13").Value = 5001
$ws.Range("L13").Value = 1000
$ws.Range("M13").Value = -4861
$ws.Range("N13").Value = -1280
$ws.Range("H14").Value = 1252.5
$ws.Range("J14").Value = 1252.5
$ws.Range("L14").Value = 1252.5
$ws.Range("N14").Value = -1588.5
$ws.Range("H104").Value = 13490
$ws.Range("J104").Value = 13490
$ws.Range("L104").Value = 13490
$ws.Range("N104").Value = -20478
$ws.Range("H122").Value = 3000
$ws.Range("I122").Value = 3000
$ws.Range("K122").Value = 9000
$ws.Range("M122").Value = -6550
